# Fruta / hortaliza, semanal
# Insert a new weekly record at row 203 ("Vega Monumental Concepción" - Uva),
# pushing the existing rows 203:260 down to 204:261.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 203 - this shifts rows 203:260
# down to 204:261, preserving their values/formatting untouched.
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with this week's record.
$ws.Cells.Item(203, 1).Value  = 11
$ws.Cells.Item(203, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(203, 3).Value  = "Bíobío"
$ws.Cells.Item(203, 4).Value  = 45120
$ws.Cells.Item(203, 5).Value  = 8
$ws.Cells.Item(203, 6).Value  = "Fruta"
$ws.Cells.Item(203, 7).Value  = 100109
$ws.Cells.Item(203, 8).Value  = "Uva"
$ws.Cells.Item(203, 9).Value  = 100109001
$ws.Cells.Item(203, 10).Value = "Uva"
$ws.Cells.Item(203, 11).Value = "Red Globe"
$ws.Cells.Item(203, 12).Value = "Primera"
$ws.Cells.Item(203, 13).Value = 100
$ws.Cells.Item(203, 14).Value = 14000
$ws.Cells.Item(203, 15).Value = 14000
$ws.Cells.Item(203, 16).Value = 14000
$ws.Cells.Item(203, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(203, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(203, 19).Value = 1400
$ws.Cells.Item(203, 20).Value = 10
